$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-11 per regenerated save data
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 11
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 1
